$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new text value, derived from the source diff.
$updates = @{
    'D2' = '69.271.86'
    'E2' = '  -3.35%  '
    'D3' = '3.509.16'
    'E3' = '  -4.87%  '
    'E4' = '  -0.13%  '
    'D5' = '582.37'
    'E5' = '  -1.21%  '
    'D6' = '173.95'
    'E6' = '  -3.92%  '
    'E7' = '  +0.77%  '
    'D8' = '3.500.42'
    'E8' = '  -4.96%  '
    'E9' = '  +0.02%  '
    'D10' = '0.189'
    'E10' = '  -6.37%  '
    'D11' = '6.70'
    'E11' = '  +4.95%  '
    'D12' = '0.595'
    'E12' = '  -3.29%  '
    'D13' = '46.98'
    'E13' = '  -6.18%  '
    'D14' = '0.0000277'
    'E14' = '  -3.99%  '
    'D15' = '674.81'
    'E15' = '  -1.61%  '
    'D16' = '4.067.39'
    'E16' = '  -5.15%  '
    'E17' = '  -3.82%  '
    'D18' = '69.171.28'
    'E18' = '  -3.60%  '
    'D19' = '3.505.58'
    'E19' = '  -4.95%  '
    'E20' = '  -1.35%  '
    'E21' = '  -4.19%  '
    'D22' = '11.18'
    'E22' = '  -4.48%  '
    'D23' = '0.903'
    'E23' = '  -4.56%  '
    'D24' = '16.10'
    'E24' = '  -10.01%  '
    'D25' = '97.89'
    'E25' = '  -5.90%  '
    'E26' = '  -4.45%  '
    'E27' = '  -0.57%  '
    'E28' = '  +0.08%  '
    'D29' = '2.65'
    'E29' = '  -7.03%  '
    'D30' = '9.44'
    'E30' = '  -7.62%  '
    'D31' = '32.96'
    'E31' = '  -6.95%  '
    'D32' = '8.71'
    'E32' = '  -6.23%  '
    'D33' = '3.18'
    'E33' = '  -8.37%  '
    'B34' = 'Mantle'
    'C34' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D34' = '1.36'
    'E34' = '  -6.05%  '
    'B35' = 'NEARProtocol'
    'C35' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D35' = '7.28'
    'E35' = '  -1.39%  '
    'D36' = '594.62'
    'E36' = '  +4.98%  '
    'D37' = '3.61'
    'E37' = '  -15.56%  '
    'D38' = '10.88'
    'E38' = '  -3.93%  '
    'D39' = '0.104'
    'E39' = '  -5.13%  '
    'D40' = '57.33'
    'E40' = '  -3.73%  '
    'D41' = '1.00'
    'E41' = '  +0.16%  '
    'D42' = '0.0438'
    'E42' = '  -6.34%  '
    'D43' = '0.335'
    'E43' = '  -5.32%  '
    'E44' = '  -6.91%  '
    'D45' = '3.411.23'
    'E45' = '  -9.42%  '
    'D46' = '33.36'
    'E46' = '  -6.46%  '
    'D47' = '0.0₃0708'
    'E47' = '  -9.13%  '
    'E48' = '  -0.71%  '
    'D49' = '2.60'
    'E49' = '  -7.66%  '
    'E50' = '  -0.85%  '
    'D51' = '5.79'
    'E51' = '  +18.27%  '
}

# These cells hold plain-text values (prices / percentages) that look numeric
# (e.g. "173.95", "1.00"). Excel would otherwise silently convert them to real
# numbers when assigned through .Value, so force Text format on the whole
# Price/Volume column range first, then restore the default "Normal" style
# afterwards so no extra cell-level formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$dataRange.Style = "Normal"

